$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# NgaDV daily report update: append a new day block (rows 40-48) describing
# 13/10/2014 (value 41925), mirroring the existing day blocks above it.
#
# Style bookkeeping: we need three brand-new cellXfs entries, created in this
# exact order so they land at indices 10, 11 and 12 (matching the target
# OOXML):
#   10 -> wrapText only                         (used by C41, C42)
#   11 -> numFmtId 14 (date) + fill "00B0F0"     (used by B40)
#   12 -> fill "FFC000" only                     (used by B41..B48)
# ---------------------------------------------------------------------------

# --- Step 1: create style #10 (wrapText) on C41 ----------------------------
$ws.Range("C41").WrapText = $true

# --- Step 2: create style #11 (date fmt + blue fill) on B40 ----------------
$ws.Range("B40").Value = 41925
$ws.Range("B40").Interior.Color = 15773696   # RGB(0,176,240) == fill "FF00B0F0"
$ws.Range("B40").NumberFormat = "mm-dd-yy"   # maps to built-in numFmtId 14

# C40 gets the same format as the other date rows' adjoining cell (fill only)
$ws.Range("C4").Copy()
$ws.Range("C40").PasteSpecial(-4122)
$ws.Range("C40").ClearContents()

# --- Step 3: create style #12 (orange fill only) on B41 --------------------
$ws.Range("B41").Interior.Color = 49407      # RGB(255,192,0) == fill "FFFFC000"

# Apply the same style (#12) to the rest of column B in this block
$ws.Range("B41").Copy()
$ws.Range("B42:B48").PasteSpecial(-4122)

# C42 also needs the wrapText style (#10), same as C41
$ws.Range("C41").Copy()
$ws.Range("C42").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Values
# ---------------------------------------------------------------------------

# Row 40: date header (value already set above)

# Row 41
$ws.Range("B41").Value = "Kế hoạch"
$ws.Range("C41").Value = "- Công việc 1:Start coding, thêm appcompat cho sicco app, Xuống version android, `n- Công việc 2: "

# Row 42
$ws.Range("B42").Value = "Kết quả đạt được"
$ws.Range("C42").Value = "- Công việc 1: 50%`n- Công việc 2: "

# Row 43
$ws.Range("B43").Value = "Trạng thái"
$ws.Range("C43").Value = "Chậm tiến độ"

# Row 44
$ws.Range("B44").Value = "Vấn đề gặp phải"
$ws.Range("C44").Value = "- Vấn đề 1: không down được android 4.0.3 , do wifi laptop có vấn đề."

# Row 45 (B45 left blank, merges visually with B44 in the source layout)
$ws.Range("C45").Value = "- Vấn đề 2:"

# Row 46
$ws.Range("B46").Value = "Giải quyết vấn đề:"
$ws.Range("C46").Value = "- Vấn đề 1: sửa lại laptop."

# Row 47 (B47 left blank)
$ws.Range("C47").Value = "- Vấn đề 2: giải quyết như sau...."

# Row 48
$ws.Range("B48").Value = "Kế hoạch ngày mai"
$ws.Range("C48").Value = "Bắt đầu code buid giao diện cho các màn hình liên quan đến Công văn ,Công việc."

# The two wrapped cells (C41, C42) must keep the sheet's default row height
# -- re-fit them so no stray ht="..." customHeight="1" sticks around.
$ws.Rows("41:42").AutoFit()

# ---------------------------------------------------------------------------
# View: mirror the saved sheetView (scrolled to show the new rows, with B40
# selected)
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("B40").Select()
